# Apply the "update test documents and update Language page" edit:
#  - Add two new API rows (batchInsert / batchUpdate) to each of the two
#    worksheets (ILookupItemService-API and ILookupClassifyService-API),
#    modeled on the existing single-path-param "findItem" row block.
#  - Make the second worksheet the active tab / selected sheet, with a
#    new selection over the freshly added rows; sheet1 loses its old
#    "active" scroll/selection state.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # ILookupItemService-API
$ws2 = $wb.Worksheets.Item(2)   # ILookupClassifyService-API

# ---------------------------------------------------------------------
# Sheet 1 (ILookupItemService-API): append batchInsert + batchUpdate,
# each a 4-row merged block copied from the existing "findItem" block
# (A9:M12), then re-label the method name / path cells (request-type
# cells are filled in afterwards, to match shared-string insert order).
# ---------------------------------------------------------------------

$ws1.Range("A9:M12").Copy($ws1.Range("A13:M16"))
$ws1.Range("A13").Value = "batchInsert"
$ws1.Range("C13").Value = "services/lookup/lookupItemService/batchInsert"

$ws1.Range("A9:M12").Copy($ws1.Range("A17:M20"))
$ws1.Range("A17").Value = "batchUpdate"
$ws1.Range("C17").Value = "services/lookup/lookupItemService/batchUpdate"

# ---------------------------------------------------------------------
# Sheet 2 (ILookupClassifyService-API): same two new rows, but pointing
# at the lookupClassifyService batch endpoints. The template block is
# copied from sheet 1's findItem rows too (matches the source workbook,
# including its leftover "itemId" parameter placeholder).
# ---------------------------------------------------------------------

$ws1.Range("A9:M12").Copy($ws2.Range("A15:M18"))
$ws2.Range("A15").Value = "batchInsert"
$ws2.Range("C15").Value = "services/lookup/lookupClassifyService/batchInsert"

$ws1.Range("A9:M12").Copy($ws2.Range("A19:M22"))
$ws2.Range("A19").Value = "batchUpdate"
$ws2.Range("C19").Value = "services/lookup/lookupClassifyService/batchUpdate"

# Request-type column (B) filled in last, POST before PUT.
$ws1.Range("B13").Value = "POST"
$ws2.Range("B15").Value = "POST"
$ws1.Range("B17").Value = "PUT"
$ws2.Range("B19").Value = "PUT"

# ---------------------------------------------------------------------
# View/selection state: sheet2 becomes the active tab with the new rows
# selected; sheet1 keeps a plain selection over its new rows too.
# ---------------------------------------------------------------------

$ws1.Activate()
$ws1.Range("B13:B20").Select()

$ws2.Activate()
$ws2.Range("B15:B22").Select()

Write-Host "Edit complete"
